# registroEntrada.xlsx — add 3 new "entrada" records (rows 9-11) and correct
# the existing record in row 8 (code, destination and timestamp).
#
# Excel auto-coerces plain numeric-looking strings ("15", "1", DNI numbers, …)
# into real numbers when assigned through Range.Value, which would store them
# as <c t="n"> instead of the shared-string <c t="s"> cells used throughout
# this sheet. To keep every cell a text value (matching the original
# workbook's convention) we stage the literal strings in a scratch area that
# has been formatted as Text ("@") first, copy that staged block, and paste
# *values only* into the destination cells — this leaves the destination
# cells' existing border/style (s="3") completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Correct the existing record in row 8
#    Código 5 -> 15, destino/fecha updated to the corrected trip details
# ---------------------------------------------------------------------
$ws.Range("F8").Value = "Áncash - Aija - Coris"
$ws.Range("G8").Value = "2023-12-05 09:11:20"

$scratch = $ws.Range("Z500")
$scratch.NumberFormat = "@"
$scratch.Value = "15"
$scratch.Copy()
$ws.Range("A8").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0
$scratch.Clear()

# ---------------------------------------------------------------------
# 2. Append three new records as rows 9, 10 and 11
# ---------------------------------------------------------------------

# 2a. Give the new rows the same border/style as row 8 first.
$ws.Range("A8:J8").Copy()
$ws.Range("A9:J11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# 2b. Stage the new rows' text in a scratch block formatted as Text, then
#     paste only the values into the real destination range.
$newRows = $ws.Range("A200:J202")
$newRows.NumberFormat = "@"

$record9  = @("16","98765687","isaac iva takamura rojas","BXI-IJT","Bus de Transporte","Amazonas - Luya - Inguilpata","2023-12-05 09:27:15","Lorenzo Daniel Arteaga Gordillo","900.0","1")
$record10 = @("17","75156514","Lorenzo Daniel Arteaga Gordillo","B89-PIJ","Vehiculo del Personal","Áncash - Santa - Santa","2023-12-05 14:46:16","Lorenzo Daniel Arteaga Gordillo","100.0","1")
$record11 = @("18","75156514","Lorenzo Daniel Arteaga Gordillo","BXI-IJT","Bus de Transporte","Apurímac - Abancay - Abancay","2023-12-07 20:14:11","Lorenzo Daniel Arteaga Gordillo","150.0","1")

$records = @($record9, $record10, $record11)

for ($r = 0; $r -lt 3; $r++) {
  for ($c = 0; $c -lt 10; $c++) {
    $newRows.Cells.Item($r + 1, $c + 1).Value = $records[$r][$c]
  }
}

$newRows.Copy()
$ws.Range("A9:J11").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0
$newRows.Clear()

# ---------------------------------------------------------------------
# 3. Column D/E best-fit widths widen slightly to accommodate the new
#    (longer) values — approximate with AutoFit.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

Write-Host "Applied registroEntrada corrections: row 8 updated, rows 9-11 added."
